$d = $word.ActiveDocument

# 1) Rename the field id: mobilizedResource -> resourceInfo
$r1 = $d.Content.Find.Execute("mobilizedResource", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "resourceInfo", 2)

# 2) Shorten the French label
$r2 = $d.Content.Find.Execute("Ressource engagée / à engager", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "Ressource", 2)

# 3) Collapse the two-sentence description (joined by a manual line break) into a single sentence
$find3 = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés en 15-15 et 15-SMUR pour le message RS-RI^lObjet permettant de communiquer la liste des ressources à engager en 15-SMUR pour le message RS-ER"
$replace3 = "Objet permettant de communquer la liste des ressource et vecteurs mobilisés en 15-15 et 15-SMUR"
$r3 = $d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $replace3, 2)

Write-Host "mobilizedResource->resourceInfo: $r1"
Write-Host "Ressource label: $r2"
Write-Host "Description merge: $r3"
